# Kpi_Item_Report.xlsx — "thêm đơn hàng trực tiếp vào báo cáo kpi"
#
# The sheet already has a "Gián tiếp" (Indirect) 12-column block (D:O) made of
# 4 sub-groups of 3 columns each (Kế hoạch / Thực hiện / Tỉ lệ (%)):
#   D:F  Sản lượng theo đơn gián tiếp
#   G:I  Doanh số theo đơn hàng gián tiếp
#   J:L  Số đơn hàng gián tiếp
#   M:O  Số đại lý theo đơn gián tiếp
#
# This change appends an identical "Trực tiếp" (Direct) 12-column block in
# P:AA, mirroring the same layout/format:
#   P:R  Sản lượng theo đơn trực tiếp
#   S:U  Doanh số theo đơn hàng trực tiếp
#   V:X  Số đơn hàng trực tiếp
#   Y:AA Số đại lý theo đơn trực tiếp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clone the formatting + values of the whole Indirect block (rows 7-10,
#    D:O) onto the new Direct block (rows 7-10, P:AA), two passes per range
#    (xlPasteFormats then xlPasteValues) so every font/fill/border/number
#    format/alignment matches exactly without creating duplicate styles.
#    The placeholder/label text is overwritten with the Direct-specific
#    strings right afterwards.
# ---------------------------------------------------------------------
$ws.Range("D7:O8").Copy()
$ws.Range("P7:AA8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D7:O8").Copy()
$ws.Range("P7:AA8").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("B9").Copy()
$ws.Range("P9:Z9").PasteSpecial(-4122)
$ws.Range("O9").Copy()
$ws.Range("AA9").PasteSpecial(-4122)

$ws.Range("D10:O10").Copy()
$ws.Range("P10:AA10").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D10:O10").Copy()
$ws.Range("P10:AA10").PasteSpecial(-4163) # xlPasteValues

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Group headers (row 7) — one label per 3-column sub-group.
# ---------------------------------------------------------------------
$ws.Range("P7").Value = "Sản lượng theo đơn trực tiếp"
$ws.Range("S7").Value = "Doanh số theo đơn hàng trực tiếp"
$ws.Range("V7").Value = "Số đơn hàng trực tiếp"
$ws.Range("Y7").Value = "Số đại lý theo đơn trực tiếp"

# ---------------------------------------------------------------------
# 3. Merge the new group headers exactly like the existing ones
#    (D7:F7, G7:I7, J7:L7, M7:O7).
# ---------------------------------------------------------------------
$ws.Range("P7:R7").Merge()
$ws.Range("S7:U7").Merge()
$ws.Range("V7:X7").Merge()
$ws.Range("Y7:AA7").Merge()

# Extend the "template-row" merge (was A9:O9) to cover the new columns.
$ws.Range("A9:O9").UnMerge()
$ws.Range("A9:AA9").Merge()

# ---------------------------------------------------------------------
# 4. Row 8 sub-headers (Kế hoạch / Thực hiện / Tỉ lệ (%)) already carry the
#    right text because they were cloned from D8:O8, which already cycles
#    through those same three labels — nothing further to set there.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5. Row 10 placeholders for the new Direct columns.
# ---------------------------------------------------------------------
$ws.Range("P10").Value = "{{KpiItemReports.Lines.DirectQuantityPlanned}}"
$ws.Range("Q10").Value = "{{KpiItemReports.Lines.DirectQuantity}}"
$ws.Range("R10").Value = "{{KpiItemReports.Lines.DirectQuantityRatio}}"
$ws.Range("S10").Value = "{{KpiItemReports.Lines.DirectRevenuePlanned}}"
$ws.Range("T10").Value = "{{KpiItemReports.Lines.DirectRevenue}}"
$ws.Range("U10").Value = "{{KpiItemReports.Lines.DirectRevenueRatio}}"
$ws.Range("V10").Value = "{{KpiItemReports.Lines.DirectAmountPlanned}}"
$ws.Range("W10").Value = "{{KpiItemReports.Lines.DirectAmount}}"
$ws.Range("X10").Value = "{{KpiItemReports.Lines.DirectAmountRatio}}"
$ws.Range("Y10").Value = "{{KpiItemReports.Lines.DirectStorePlanned}}"
$ws.Range("Z10").Value = "{{KpiItemReports.Lines.DirectStore}}"
$ws.Range("AA10").Value = "{{KpiItemReports.Lines.DirectStoreRatio}}"

# ---------------------------------------------------------------------
# 6. Keep the printable/used range consistent with the new columns and
#    mirror the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("G17").Select()
